$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new C (nombre_aides), new E (montant_total)
$updates = @{
    9   = @{ C = 69574;  E = 191367649 }
    125 = @{ C = 4597;   E = 13136795 }
    150 = @{ C = 95010;  E = 278992606 }
    167 = @{ C = 101529; E = 194966146 }
    168 = @{ C = 285008; E = 1210390222 }
    169 = @{ C = 562600; E = 1284501384 }
    170 = @{ C = 367381; E = 2845626697 }
    171 = @{ C = 115156; E = 446646135 }
    173 = @{ C = 54389;  E = 151875013 }
    174 = @{ C = 357235; E = 1017861449 }
    175 = @{ C = 125553; E = 812943051 }
    177 = @{ C = 96756;  E = 174740928 }
    178 = @{ C = 75360;  E = 102747363 }
    179 = @{ C = 235713; E = 812627643 }
    180 = @{ C = 141482; E = 340835117 }
    181 = @{ C = 7894;   E = 11291524 }
    188 = @{ C = 19707;  E = 66147299 }
    204 = @{ C = 4759;   E = 11763349 }
    259 = @{ C = 6240;   E = 14751339 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("E$row").Value = $vals.E
}

$wb.Save()
